# ---------------------------------------------------------------------------
# empleados.xlsx: add "localidad" (split out of "direccion") and
# "horario_salida" columns, and load the remaining employee rows.
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1): insert "localidad" after "direccion" and
# "horario_salida" after "horario_ingreso" -- pushes preferencia/observaciones
# to columns G/H.
# ---------------------------------------------------------------------------
$ws.Cells.Item(1,1).Value = "id_empleado"
$ws.Cells.Item(1,2).Value = "nombre"
$ws.Cells.Item(1,3).Value = "direccion"
$ws.Cells.Item(1,4).Value = "localidad"
$ws.Cells.Item(1,5).Value = "horario_ingreso"
$ws.Cells.Item(1,6).Value = "horario_salida"
$ws.Cells.Item(1,7).Value = "preferencia"
$ws.Cells.Item(1,8).Value = "observaciones"

# ---------------------------------------------------------------------------
# Data rows 2-8
# ---------------------------------------------------------------------------
# Row 2 - Ana Perez (address split into direccion/localidad)
$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = "Ana Pérez"
$ws.Cells.Item(2,3).Value = "Av Corrientes 1234"
$ws.Cells.Item(2,4).Value = "CABA"
$ws.Cells.Item(2,5).Value = 8
$ws.Cells.Item(2,6).Value = 16
$ws.Cells.Item(2,7).Value = "ninguna"

# Row 3 - Cristian
$ws.Cells.Item(3,1).Value = 2
$ws.Cells.Item(3,2).Value = "Cristian"
$ws.Cells.Item(3,3).Value = "emilio mitre 772"
$ws.Cells.Item(3,4).Value = "CABA"
$ws.Cells.Item(3,5).Value = 9
$ws.Cells.Item(3,6).Value = 17

# Row 4 - juan
$ws.Cells.Item(4,1).Value = 3
$ws.Cells.Item(4,2).Value = "juan"
$ws.Cells.Item(4,3).Value = "tomas valle 772"
$ws.Cells.Item(4,4).Value = "La Matanza"
$ws.Cells.Item(4,5).Value = 8
$ws.Cells.Item(4,6).Value = 17

# Row 5 - raul
$ws.Cells.Item(5,1).Value = 4
$ws.Cells.Item(5,2).Value = "raul"
$ws.Cells.Item(5,3).Value = "Corrientes 2445"
$ws.Cells.Item(5,4).Value = "caba"
$ws.Cells.Item(5,5).Value = 8
$ws.Cells.Item(5,6).Value = 16

# Row 6 - Mica
$ws.Cells.Item(6,1).Value = 5
$ws.Cells.Item(6,2).Value = "Mica"
$ws.Cells.Item(6,3).Value = "Av Militar 3900"
$ws.Cells.Item(6,4).Value = "ciudadela"
$ws.Cells.Item(6,5).Value = 8
$ws.Cells.Item(6,6).Value = 17

# Row 7 - rosa
$ws.Cells.Item(7,1).Value = 6
$ws.Cells.Item(7,2).Value = "rosa"
$ws.Cells.Item(7,3).Value = "Benjamín Franklin 2298"
$ws.Cells.Item(7,4).Value = "paso del rey"
$ws.Cells.Item(7,5).Value = 9
$ws.Cells.Item(7,6).Value = 17

# Row 8 - jorge
$ws.Cells.Item(8,1).Value = 7
$ws.Cells.Item(8,2).Value = "jorge"
$ws.Cells.Item(8,3).Value = "Roseti 1601"
$ws.Cells.Item(8,4).Value = "caba"
$ws.Cells.Item(8,5).Value = 9
$ws.Cells.Item(8,6).Value = 16

# ---------------------------------------------------------------------------
# Formatting
# ---------------------------------------------------------------------------
# Header: bold, centered, wrap text (matches existing header style)
$header = $ws.Range("A1:H1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108   # xlCenter
$header.WrapText = $true

# Leftover per-cell wrap-text formatting (mirrors the source file's uneven
# formatting carried over from copy/pasted rows)
$ws.Range("B2").WrapText = $true
$ws.Range("D2").WrapText = $true
$ws.Range("G2").WrapText = $true
$ws.Range("C6:C8").WrapText = $true

# Column C (direccion) width (~15.3 chars)
$ws.Columns.Item(3).ColumnWidth = 14.5

# Leftover empty-but-formatted cells (underline), matching source formatting
$ws.Range("G3").Font.Underline = $true
$ws.Range("H6").Font.Underline = $true

# Row heights (row 2 keeps its pre-existing explicit height; row 7 grew to a
# second line like the header row)
$ws.Rows.Item(2).RowHeight = 15.75
$ws.Rows.Item(7).RowHeight = 23.85

$ws.Range("H15").Select() | Out-Null
